# The deck ships with two embedded theme parts:
#   - the "Integral" theme, currently applied to the presentation's only
#     slide master (and therefore to every slide), and
#   - a spare/leftover "Office Theme" (used only by the notes master).
#
# The author switched the presentation's applied Design from "Integral"
# to "Office Theme" (Design tab -> Office Theme), which recolors the
# slide master - and every slide built on it - with the Office Theme
# palette. Re-apply that by writing the Office Theme's 12 theme colors
# into the color scheme of the presentation's active design/master.

function Convert-HexToComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Office Theme color scheme (standard PowerPoint "Office Theme" design).
# Index order matches the OOXML clrScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeTheme = @{
    1  = "000000"
    2  = "FFFFFF"
    3  = "44546A"
    4  = "E7E6E6"
    5  = "5B9BD5"
    6  = "ED7D31"
    7  = "A5A5A5"
    8  = "FFC000"
    9  = "4472C4"
    10 = "70AD47"
    11 = "0563C1"
    12 = "954F72"
}

foreach ($idx in $officeTheme.Keys) {
    $colors.Item($idx).RGB = Convert-HexToComRGB $officeTheme[$idx]
}

Write-Output "Applied Office Theme colors to the presentation design."
Write-Output "Design name: $($design.Name)"
